$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.282.18"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.864.54"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2872"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06560"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07881"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").Value = "1.865.48"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6936"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.094"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "265.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "30.255.65"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007653"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("D21").Value = "2.115.01"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.232"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.372"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.939"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.354"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09869"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.457"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.060"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04754"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7008"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.795"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.202"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8424"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.117"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "937.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.082"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05677"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.32%  "
